$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Value" column (D) for resistors R1 (row 8), R2 (row 9), R4 (row 11)
$ws.Range("D8").Value = "2k"
$ws.Range("D9").Value = "221R"
$ws.Range("D11").Value = "10K"

# Update the active selection on the sheet
$ws.Range("A3:I17").Select()

$wb.Save()
